$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 0: sheet view (split pane scrolled to the new data + selection) ---
$ws.Activate()
$ws.Range("L1").Select()
$excel.ActiveWindow.SplitRow = 92
$excel.ActiveWindow.SplitColumn = 0
$excel.ActiveWindow.Split = $true
$ws.Range("D114").Select()

# --- Step 1: rows 82-91, fill previously-empty B/K/L cells with values and clear their styles ---
$ws.Range("B82").ClearFormats()
$ws.Range("B82").Value = 1.51563076473506
$ws.Range("K82:L82").ClearFormats()
$ws.Range("K82").Value = 1.3792365585475499
$ws.Range("L82").Value = 1.4179243953946901
$ws.Range("B83").ClearFormats()
$ws.Range("B83").Value = 1.54257973490762
$ws.Range("K83:L83").ClearFormats()
$ws.Range("K83").Value = 1.45323133960696
$ws.Range("L83").Value = 1.45441235336025
$ws.Range("B84").ClearFormats()
$ws.Range("B84").Value = 1.5975597012243801
$ws.Range("K84:L84").ClearFormats()
$ws.Range("K84").Value = 1.5022230591250401
$ws.Range("L84").Value = 1.50631529067879
$ws.Range("B85").ClearFormats()
$ws.Range("B85").Value = 1.5423757127992099
$ws.Range("K85:L85").ClearFormats()
$ws.Range("K85").Value = 1.39086458261685
$ws.Range("L85").Value = 1.48384425703649
$ws.Range("B86").ClearFormats()
$ws.Range("B86").Value = 1.55771011715065
$ws.Range("K86:L86").ClearFormats()
$ws.Range("K86").Value = 1.46374552789837
$ws.Range("L86").Value = 1.46518448598394
$ws.Range("B87").ClearFormats()
$ws.Range("B87").Value = 1.48422449210259
$ws.Range("K87:L87").ClearFormats()
$ws.Range("K87").Value = 1.3856829643775099
$ws.Range("L87").Value = 1.3866860755928401
$ws.Range("B88").ClearFormats()
$ws.Range("B88").Value = 1.5975001003224301
$ws.Range("K88:L88").ClearFormats()
$ws.Range("K88").Value = 1.55259740450018
$ws.Range("L88").Value = 1.5215908670937
$ws.Range("B89").ClearFormats()
$ws.Range("B89").Value = 1.5706950106781701
$ws.Range("K89:L89").ClearFormats()
$ws.Range("K89").Value = 1.49121056249448
$ws.Range("L89").Value = 1.4802731985691899
$ws.Range("B90").ClearFormats()
$ws.Range("B90").Value = 1.54559117409575
$ws.Range("K90:L90").ClearFormats()
$ws.Range("K90").Value = 1.47611314589978
$ws.Range("L90").Value = 1.4428712684556499
$ws.Range("B91").ClearFormats()
$ws.Range("B91").Value = 1.4961218538771399
$ws.Range("K91:L91").ClearFormats()
$ws.Range("K91").Value = 1.3980895490257701
$ws.Range("L91").Value = 1.4154019259042501

# --- Step 2: rows 92-101, write A:N values via array assignment ---
$row92 = New-Object 'object[,]' 1,14
$row92[0,0] = 10
$row92[0,1] = 1.61934427605201
$row92[0,2] = 1.54324248058308
$row92[0,3] = 1.530604286952
$row92[0,4] = 1.5279060597079299
$row92[0,5] = 1.4806564985478801
$row92[0,6] = 1.5074008326771999
$row92[0,7] = 1.48495677057885
$row92[0,8] = 1.4907491758035101
$row92[0,9] = 1.54228611681689
$row92[0,10] = 1.5914524054152599
$row92[0,11] = 1.58870257241544
$row92[0,12] = 1.48464533321087
$row92[0,13] = 1.5127293540110101
$ws.Range("A92:N92").Value = $row92
$row93 = New-Object 'object[,]' 1,14
$row93[0,0] = 10
$row93[0,1] = 1.7000773394984501
$row93[0,2] = 1.79846179533265
$row93[0,3] = 1.75447451055571
$row93[0,4] = 1.81338215141631
$row93[0,5] = 1.7515986871024001
$row93[0,6] = 1.7998118956960301
$row93[0,7] = 1.7610410289402201
$row93[0,8] = 1.7690833562595101
$row93[0,9] = 1.5237035273520301
$row93[0,10] = 1.70813460982066
$row93[0,11] = 1.7414208710678001
$row93[0,12] = 1.7726806321142801
$row93[0,13] = 1.78500891714931
$ws.Range("A93:N93").Value = $row93
$row94 = New-Object 'object[,]' 1,14
$row94[0,0] = 10
$row94[0,1] = 1.85389438797039
$row94[0,2] = 1.8659030146579301
$row94[0,3] = 1.8414617466840699
$row94[0,4] = 1.8789898289781
$row94[0,5] = 1.8319551385241399
$row94[0,6] = 1.8718506777848001
$row94[0,7] = 1.84451940067207
$row94[0,8] = 1.8428661100552
$row94[0,9] = 1.6632789158230601
$row94[0,10] = 1.8434675684816999
$row94[0,11] = 1.85532648952408
$row94[0,12] = 1.8492203230240301
$row94[0,13] = 1.85958759031402
$ws.Range("A94:N94").Value = $row94
$row95 = New-Object 'object[,]' 1,14
$row95[0,0] = 10
$row95[0,1] = 1.7659163423676001
$row95[0,2] = 1.79708368378766
$row95[0,3] = 1.74359296749018
$row95[0,4] = 1.81195262961022
$row95[0,5] = 1.7291172813192599
$row95[0,6] = 1.7855088554976499
$row95[0,7] = 1.7475663000766699
$row95[0,8] = 1.76044111963758
$row95[0,9] = 1.5405779640111299
$row95[0,10] = 1.7567549217584499
$row95[0,11] = 1.7628442895315399
$row95[0,12] = 1.76167439020335
$row95[0,13] = 1.7792347130527599
$ws.Range("A95:N95").Value = $row95
$row96 = New-Object 'object[,]' 1,14
$row96[0,0] = 10
$row96[0,1] = 1.38982775475084
$row96[0,2] = 1.2278879198974999
$row96[0,3] = 1.2267744510871199
$row96[0,4] = 1.1804314880223501
$row96[0,5] = 1.19752057434684
$row96[0,6] = 1.19525509494782
$row96[0,7] = 1.42557572682949
$row96[0,8] = 1.3615345059266399
$row96[0,9] = 1.34654481417787
$row96[0,10] = 1.39985558280322
$row96[0,11] = 1.3672492210529099
$row96[0,12] = 1.3579153127407599
$row96[0,13] = 1.37334920071543
$ws.Range("A96:N96").Value = $row96
$row97 = New-Object 'object[,]' 1,14
$row97[0,0] = 10
$row97[0,1] = 1.49016462345485
$row97[0,2] = 1.6524052009576999
$row97[0,3] = 1.3696100871219199
$row97[0,4] = 1.35315082583033
$row97[0,5] = 1.31907575034854
$row97[0,6] = 1.2981208490398
$row97[0,7] = 1.30489711376144
$row97[0,8] = 1.26122990273998
$row97[0,9] = 1.53878500185916
$row97[0,10] = 1.44111093543391
$row97[0,11] = 1.37684480439765
$row97[0,12] = 1.3033661450552601
$row97[0,13] = 1.30314441670938
$ws.Range("A97:N97").Value = $row97
$row98 = New-Object 'object[,]' 1,14
$row98[0,0] = 10
$row98[0,1] = 1.46855143861149
$row98[0,2] = 1.8847556661648599
$row98[0,3] = 1.3642747606743
$row98[0,4] = 1.37018198077834
$row98[0,5] = 1.3051689189834701
$row98[0,6] = 1.2893313364648999
$row98[0,7] = 1.32983725828248
$row98[0,8] = 1.2492532601386801
$row98[0,9] = 1.54468143178101
$row98[0,10] = 1.40732322679189
$row98[0,11] = 1.4157107718312301
$row98[0,12] = 1.2935825509698
$row98[0,13] = 1.33095374159687
$ws.Range("A98:N98").Value = $row98
$row99 = New-Object 'object[,]' 1,14
$row99[0,0] = 10
$row99[0,1] = 1.55845610118504
$row99[0,2] = 1.76133126255459
$row99[0,3] = 1.42540943173039
$row99[0,4] = 1.3809044695706101
$row99[0,5] = 1.3226799053488101
$row99[0,6] = 1.2527819447064401
$row99[0,7] = 1.3394697193601299
$row99[0,8] = 1.2081581745124399
$row99[0,9] = 1.64107593297887
$row99[0,10] = 1.5086797983776099
$row99[0,11] = 1.45415746047362
$row99[0,12] = 1.28060738450814
$row99[0,13] = 1.33987074322099
$ws.Range("A99:N99").Value = $row99
$row100 = New-Object 'object[,]' 1,14
$row100[0,0] = 10
$row100[0,1] = 1.54818332123055
$row100[0,2] = 1.6156725757172901
$row100[0,3] = 1.4828014587573499
$row100[0,4] = 1.4041963737188401
$row100[0,5] = 1.42687067086681
$row100[0,6] = 1.3771800832608201
$row100[0,7] = 1.3198240596738799
$row100[0,8] = 1.26643146021018
$row100[0,9] = 1.5548664035303399
$row100[0,10] = 1.5226773096471999
$row100[0,11] = 1.4851345847766799
$row100[0,12] = 1.3215435285039101
$row100[0,13] = 1.32872350681587
$ws.Range("A100:N100").Value = $row100
$row101 = New-Object 'object[,]' 1,14
$row101[0,0] = 10
$row101[0,1] = 1.4024104084150899
$row101[0,2] = 1.84937912495301
$row101[0,3] = 1.3526895215336401
$row101[0,4] = 1.3451199875044499
$row101[0,5] = 1.2582126681054
$row101[0,6] = 1.4305393552119701
$row101[0,7] = 1.3546475201261801
$row101[0,8] = 1.28647559322669
$row101[0,9] = 1.3622297874731
$row101[0,10] = 1.3490208347116199
$row101[0,11] = 1.38837871970018
$row101[0,12] = 1.3316671036696399
$row101[0,13] = 1.3346609585795499
$ws.Range("A101:N101").Value = $row101

# --- Step 3: rows 92-101, write O:S values via array assignment ---
$colO = New-Object 'object[,]' 10,1
$colO[0,0] = 1
$colO[1,0] = 2
$colO[2,0] = 3
$colO[3,0] = 4
$colO[4,0] = 5
$colO[5,0] = 6
$colO[6,0] = 7
$colO[7,0] = 8
$colO[8,0] = 9
$colO[9,0] = 10
$ws.Range("O92:O101").Value = $colO
$colP = New-Object 'object[,]' 10,1
$colP[0,0] = 113
$colP[1,0] = 117
$colP[2,0] = 118
$colP[3,0] = 131
$colP[4,0] = 422
$colP[5,0] = 606
$colP[6,0] = 607
$colP[7,0] = 608
$colP[8,0] = 609
$colP[9,0] = 610
$ws.Range("P92:P101").Value = $colP
$colQ = New-Object 'object[,]' 10,1
$colQ[0,0] = 73
$colQ[1,0] = 73
$colQ[2,0] = 73
$colQ[3,0] = 73
$colQ[4,0] = 73
$colQ[5,0] = 73
$colQ[6,0] = 73
$colQ[7,0] = 73
$colQ[8,0] = 73
$colQ[9,0] = 65
$ws.Range("Q92:Q101").Value = $colQ
$colR = New-Object 'object[,]' 10,1
$colR[0,0] = 96
$colR[1,0] = 96
$colR[2,0] = 96
$colR[3,0] = 96
$colR[4,0] = 96
$colR[5,0] = 96
$colR[6,0] = 96
$colR[7,0] = 96
$colR[8,0] = 96
$colR[9,0] = 96
$ws.Range("R92:R101").Value = $colR
$colS = New-Object 'object[,]' 10,1
$colS[0,0] = 29
$colS[1,0] = 29
$colS[2,0] = 29
$colS[3,0] = 29
$colS[4,0] = 29
$colS[5,0] = 29
$colS[6,0] = 29
$colS[7,0] = 29
$colS[8,0] = 29
$colS[9,0] = 29
$ws.Range("S92:S101").Value = $colS

# --- Step 4: apply cell formats (style indices) by copying from donor cells with the same style ---
$ws.Range("O2").Copy() | Out-Null
$ws.Range("O92:O101").PasteSpecial(-4122) | Out-Null
$ws.Range("P12").Copy() | Out-Null
$ws.Range("P92:P101").PasteSpecial(-4122) | Out-Null
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("Q92:Q101").PasteSpecial(-4122) | Out-Null
$ws.Range("S72").Copy() | Out-Null
$ws.Range("R92:S101").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# style 23 (theme fill, no alignment) is not pre-applied to any individual cell in the source file;
# manufacture a donor via a scratch cell far away in column S, whose <col> default style is 23
$ws.Range("S500").Value = 1
$ws.Range("S500").Copy() | Out-Null
$ws.Range("K94:N94").PasteSpecial(-4122) | Out-Null
$ws.Range("B95:J95").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("S500").Clear()
$ws.Range("D114").Select()